# Atomix perf updates — add two "new API" measurement columns to the
# "Atomix" sheet, shifting the old "SC-SC-OPT" / "unsound ts_get()..."
# columns two places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atomix")

# Insert two blank columns at D:E — this shifts the existing D/E
# (SC-SC-OPT / unsound ts_get() columns) to F/G.
$ws.Columns("D:E").Insert()

# New header labels for the inserted columns (the inserted blank columns
# already carry the same header style "s=18" as the rest of row 1).
$ws.Range("D1").Value = "new API"
$ws.Range("E1").Value = "new API (no mitig copy)"

# New per-row measurements (columns D and E).
$newApi = @{
    2  = 73
    3  = 53
    4  = 191
    5  = 206
    6  = 185
    7  = 163
    8  = 126
    9  = 100
    10 = 196
    11 = 111
    12 = 127
    13 = 47
    14 = 47
    15 = 50
    16 = 62
    17 = 80
    18 = 92
}
$newApiNoMitig = @{
    2  = 73
    3  = 45
    4  = 185
    5  = 204
    6  = 176
    7  = 166
    8  = 127
    9  = 101
    10 = 172
    11 = 105
    12 = 123
    13 = 48
    14 = 47
    15 = 47
    16 = 58
    17 = 78
    18 = 89
}

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 4).Value = $newApi[$row]
    $ws.Cells.Item($row, 5).Value = $newApiNoMitig[$row]
}

# Column widths: D matches C's existing width, E gets its own (best-fit-ish)
# width, F/G keep the widths the old D/E columns had.
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth
$ws.Columns("E").ColumnWidth = 22.42578125
$ws.Columns("F").ColumnWidth = 13.7109375
$ws.Columns("G").ColumnWidth = 53.85546875

# Final selection, matching where editing left off.
$ws.Range("D18").Select()
